$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A63").Value = "x"
$ws.Rows.Item(63).RowHeight = 15
Write-Output "done"
